$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154, shifting existing rows 154-196 down to 155-197.
$ws.Rows.Item(154).Insert()

# Populate the new row 154 with the new record's data.
$ws.Cells.Item(154, 1).Value = 11
$ws.Cells.Item(154, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(154, 3).Value = "Bíobío"
$ws.Cells.Item(154, 4).Value = 45205
$ws.Cells.Item(154, 5).Value = 8
$ws.Cells.Item(154, 6).Value = "Fruta"
$ws.Cells.Item(154, 7).Value = 100108
$ws.Cells.Item(154, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(154, 9).Value = 100108002
$ws.Cells.Item(154, 10).Value = "Mango"
$ws.Cells.Item(154, 11).Value = "Sin especificar"
$ws.Cells.Item(154, 12).Value = "Primera"
$ws.Cells.Item(154, 13).Value = 100
$ws.Cells.Item(154, 14).Value = 10000
$ws.Cells.Item(154, 15).Value = 10000
$ws.Cells.Item(154, 16).Value = 10000
$ws.Cells.Item(154, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(154, 18).Value = "Brasil"
$ws.Cells.Item(154, 19).Value = 2500
$ws.Cells.Item(154, 20).Value = 4
